$wb = $excel.ActiveWorkbook

# --- Update selections on the existing sheets (cursor moved around while
#     building the new "Alquiler" tooling) ---
$wsVehiculos = $wb.Worksheets.Item("vehiculos")
$wsVehiculos.Range("I3").Select() | Out-Null

$wsClientes = $wb.Worksheets.Item("Clientes")
$wsClientes.Range("B2").Select() | Out-Null

$wsSolicitudes = $wb.Worksheets.Item("Solicitudes")
$wsSolicitudes.Range("C1").Select() | Out-Null

# --- Add the new "Alquiler" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Alquiler"

# --- Seed the first (and only) data row ---
$ws.Range("A1").Value = 2022121

# B1 carries a quote-prefixed text date ("05-12-2022"), matching the style
# already used for the date column on the Solicitudes sheet (numFmtId 14,
# quotePrefix). Copy that format over so no new style entry is created.
$wsSolicitudes.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Value = "'05-12-2022"

$ws.Range("C1").Value = 506060
$ws.Range("D1").Value = 22222222
$ws.Range("E1").Value = "Valery"
$ws.Range("F1").Value = "Plata"
$ws.Range("G1").Value = 150000
$ws.Range("H1").Value = 2
$ws.Range("I1").Formula = "=300000"
$ws.Range("J1").Value = 78000
$ws.Range("K1").Value = 678000

$ws.Range("L1").Select() | Out-Null
